# Regenerate save_data to use K (strikeouts) instead of Strike# in column G.
# Recomputed K values for each outing row (rows 2-32 on the active sheet).

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$kValues = @{
    2  = 0
    3  = 2
    4  = 1
    5  = 1
    6  = 1
    7  = 0
    8  = 3
    9  = 1
    10 = 2
    11 = 1
    12 = 1
    13 = 1
    14 = 0
    15 = 4
    16 = 2
    17 = 2
    18 = 5
    19 = 0
    20 = 2
    21 = 1
    22 = 1
    23 = 1
    24 = 0
    25 = 2
    26 = 0
    27 = 0
    28 = 0
    29 = 1
    30 = 2
    31 = 0
    32 = 1
}

foreach ($row in $kValues.Keys) {
    $ws.Range("G$row").Value = $kValues[$row]
}
